# "Scene 1 began, with room modelled and new assets"
# Kanban-style tasklist update: several tasks move from TODO (col A) -> DOING (col B) -> DONE (col C),
# and a brand new task is recorded as DONE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (TODO): replace the three finished/started TODO items with the next three ---
$ws.Range("A3").Value = "Build to executable outside of Unity editor"
$ws.Range("A3").Style = "Neutral"

$ws.Range("A4").Value = "Create Menu for game"
$ws.Range("A4").Style = "Neutral"

$ws.Range("A5").Value = "Find/Record audio elements"
$ws.Range("A5").Style = "Neutral"

# The old rows 6-8 no longer have a TODO (column A) entry.
$ws.Range("A6").Clear()
$ws.Range("A7").Clear()
$ws.Range("A8").Clear()

# --- Column B (DOING): shift tasks that progressed, and add newly-started work ---
$ws.Range("B6").Value = "Model first level design in Unity"
$ws.Range("B6").Style = "Neutral"

$ws.Range("B7").Value = "Create custom assets if needed"
$ws.Range("B7").Style = "Neutral"

$ws.Range("B8").Value = "Model Collision recognition in Unity"
$ws.Range("B8").Style = "Neutral"

$ws.Range("B9").Value = "Create first level design"
$ws.Range("B9").Style = "Good"

# --- Column C (DONE): newly completed tasks ---
$ws.Range("C9").Value = "Model Throwing Physics in Unity"
$ws.Range("C9").Style = "Neutral"

$ws.Range("C10").Value = "Find early assets to build with"
$ws.Range("C10").Style = "Good"

# Update the active selection to match where the user ended up working
$ws.Range("D12").Select()
